{"js": "// Replace \"in HW4\" with \"through GraphWrapper in Stage 1\" in the last\n// bullet of the ADT-changes list, so the sentence reads:\n// \"Commented out unused methods in Graph that were already tested through\n//  GraphWrapper in Stage 1, as well as all of GraphWrapperTest.java and\n//  GraphWrapper.java, for the sake of coverage.\"\n\nconst body = context.document.body;\n\nconst results = body.search(\"in HW4\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"through GraphWrapper in Stage 1\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace \"in HW4\" with \"through GraphWrapper in Stage 1\" in the last\n# bullet of the ADT-changes list, so the sentence reads:\n# \"Commented out unused methods in Graph that were already tested through\n#  GraphWrapper in Stage 1, as well as all of GraphWrapperTest.java and\n#  GraphWrapper.java, for the sake of coverage.\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$found = $find.Execute(\"in HW4\", $false, $false, $false, $false, $false, $true, 1, $false, \"through GraphWrapper in Stage 1\", 2)\n"}
